$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2) | Out-Null
}

# 1. "largerly" -> "largely"
Replace-Text "increasing largerly from 5" "increasing largely from 5"

# 2. "postion" -> "position" (Java developer paragraph)
Replace-Text "This is a postion which is also have a high demand" "This is a position which is also have a high demand"

# 3. "postions" -> "positions" (end of Java developer paragraph)
Replace-Text " postions. We believe the demand of this job will increase in the next year." " positions. We believe the demand of this job will increase in the next year."

# 4. Web developer paragraph: "keep" -> "keeps", "postion" -> "position"
Replace-Text "Web developer is a stable job. It keep a stable demand from employers. And it is in the middle postion of the ranking list. " "Web developer is a stable job. It keeps a stable demand from employers. And it is in the middle position of the ranking list. "

# 5. Software developer paragraph: "Unluckly" -> "Unluckily", add missing space, remove "the"
Replace-Text "Unluckly, it is a job which is decreasing demand from employers.In the last year, it is" "Unluckily, it is a job which is decreasing demand from employers. In last year, it is"

# 6. SQL paragraph: "Datebase" -> "Database"
Replace-Text "Datebase is the most important thing" "Database is the most important thing"

# 7. JavaScript paragraph: "databse" -> "database", "have a enough reason" -> "has enough reasons"
Replace-Text "If databse is the most important thing to store information, website must be the most important thing to show and collect information. JavaScripts as a language to catch the users" "If database is the most important thing to store information, website must be the most important thing to show and collect information. JavaScripts as a language to catch the users"
Replace-Text "input, it have a enough reason to get the second highest demand." "input, it has enough reasons to get the second highest demand."

# 8. Git paragraph: "it store the code" -> "it stores the code"
Replace-Text "it store the code" "it stores the code"

# 9. General skills intro: "skills.problem" -> "skills，problem" (fullwidth comma)
Replace-Text "communication skills.problem solving skills" "communication skills，problem solving skills"

# 10. Communication skills paragraph: "see that it own the highest" -> "see how it meets the highest"
Replace-Text "It is easy to see that it own the highest demand from employers. " "It is easy to see how it meets the highest demand from employers. "

# 11. Problem solving paragraph: "That is explain that why" -> "That is an explanation why"
Replace-Text "them. That is explain that why problem solving is an important skill." "them. That is an explanation why problem solving is an important skill."

# 13. Teamwork paragraph: "must be more easy" -> "might be easier"
Replace-Text "If we could have team work, it must be more easy to get the job." "If we could have team work, it might be easier to get the job."

# 14. 2d. Organisation skill paragraph
Replace-Text "Organisation skill is the only one generic skll which is in the top3 rank but not in our required skill set." "Organizing skill is the only one generic skill which is in the top3 rank but not in our required skill set."

# 15. Question 3 first paragraph, full rewrite
Replace-Text "We don’t have a big change of the opinion about our ideal job." "In general, the individuals’ opinions of ideal job have not changed much."

# 16. Question 3 second paragraph, full rewrite (keep bookmark in place, just change surrounding text)
Replace-Text "From the below we analysis, the ideal job of us have the high rank in the list. Even it is not all of the highest. But they are still the dream work we want to get. In addtions, who know the demand from employers? " "From what we have analyzed, the ideal job has the high rank in the list. Even it is not all of the highest. But they are still the dream work we want to get. In addition, who knows the demand from employers? "
